$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 362.17392
$ws.Range("I33").Value = 254.25
$ws.Range("J33").Value = 608.8570999999999
$ws.Range("K33").Value = 254.25
$ws.Range("L33").Value = 608.8570999999999
$ws.Range("M33").Value = -25.25
$ws.Range("N33").Value = -1066.8571
$ws.Range("H103").Value = 1392.9656
$ws.Range("I103").Value = 700.8570999999999
$ws.Range("J103").Value = 1613.1818
$ws.Range("K103").Value = 2102.5713
$ws.Range("L103").Value = 4839.5454
$ws.Range("M103").Value = -1516.5713
$ws.Range("N103").Value = -6011.5454
$ws.Range("H112").Value = 2025.6897
$ws.Range("J112").Value = 2025.6897
$ws.Range("L112").Value = 6077.0691
$ws.Range("N112").Value = -8293.069100000001
$ws.Range("H138").Value = 6738.7754
$ws.Range("J138").Value = 6854.478
$ws.Range("L138").Value = 20563.434
$ws.Range("N138").Value = -30843.434

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4930.92
$ws.Range("I2").Value = 819.9474
$ws.Range("K2").Value = 819.9474
$ws.Range("M2").Value = -706.9474
$ws.Range("H32").Value = 12469.025
$ws.Range("I32").Value = 10545.632
$ws.Range("K32").Value = 10545.632
$ws.Range("M32").Value = -10258.632
$ws.Range("H61").Value = 8046.8823
$ws.Range("I61").Value = 8379
$ws.Range("J61").Value = 7572.4287
$ws.Range("K61").Value = 8379
$ws.Range("L61").Value = 7572.4287
$ws.Range("M61").Value = -8167
$ws.Range("N61").Value = -7996.4287
$ws.Range("H63").Value = 7614.5454
$ws.Range("J63").Value = 8445
$ws.Range("L63").Value = 8445
$ws.Range("N63").Value = -9817
$ws.Range("H66").Value = 7614.5454
$ws.Range("J66").Value = 8445
$ws.Range("L66").Value = 42225
$ws.Range("N66").Value = -49089
$ws.Range("H74").Value = 2481
$ws.Range("I74").Value = 1777.4
$ws.Range("J74").Value = 3988.7144
$ws.Range("K74").Value = 1777.4
$ws.Range("L74").Value = 3988.7144
$ws.Range("M74").Value = -903.4000000000001
$ws.Range("N74").Value = -5736.7144
$ws.Range("H77").Value = 2481
$ws.Range("I77").Value = 1777.4
$ws.Range("J77").Value = 3988.7144
$ws.Range("K77").Value = 8887
$ws.Range("L77").Value = 19943.572
$ws.Range("M77").Value = -4519
$ws.Range("N77").Value = -28679.572
$ws.Range("H101").Value = 75000
$ws.Range("J101").Value = 75000
$ws.Range("L101").Value = 75000
$ws.Range("N101").Value = -81490
$ws.Range("H116").Value = 4930.92
$ws.Range("I116").Value = 819.9474
$ws.Range("K116").Value = 819.9474
$ws.Range("M116").Value = 1474.0526
$ws.Range("H122").Value = 4592.5
$ws.Range("I122").Value = 4348
$ws.Range("K122").Value = 13044
$ws.Range("M122").Value = -10594
$ws.Range("H132").Value = 4834.143
$ws.Range("I132").Value = 4347
$ws.Range("J132").Value = 9462
$ws.Range("K132").Value = 13041
$ws.Range("L132").Value = 28386
$ws.Range("M132").Value = -10511
$ws.Range("N132").Value = -33446
$ws.Range("H135").Value = 59473.816
$ws.Range("J135").Value = 59473.816
$ws.Range("L135").Value = 59473.816
$ws.Range("N135").Value = -69613.81599999999
$ws.Range("H136").Value = 8046.8823
$ws.Range("I136").Value = 8379
$ws.Range("J136").Value = 7572.4287
$ws.Range("K136").Value = 25137
$ws.Range("L136").Value = 22717.2861
$ws.Range("M136").Value = -22587
$ws.Range("N136").Value = -27817.2861
$ws.Range("H137").Value = 69798
$ws.Range("J137").Value = 69798
$ws.Range("L137").Value = 69798
$ws.Range("N137").Value = -79998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4930.92
$ws.Range("I3").Value = 819.9474
$ws.Range("K3").Value = 819.9474
$ws.Range("M3").Value = -705.9474
$ws.Range("H59").Value = 98997.5
$ws.Range("J59").Value = 115330
$ws.Range("L59").Value = 115330
$ws.Range("N59").Value = -117024
$ws.Range("H94").Value = 2901.2104
$ws.Range("I94").Value = 2901.2104
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 2901.2104
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("M94").Value = -2450.2104
$ws.Range("H105").Value = 12408.655
$ws.Range("I105").Value = 12193.056
$ws.Range("K105").Value = 12193.056
$ws.Range("M105").Value = -10446.056
$ws.Range("H134").Value = 4764.3335
$ws.Range("I134").Value = 4355.8
$ws.Range("J134").Value = 5785.6665
$ws.Range("K134").Value = 13067.4
$ws.Range("L134").Value = 17356.9995
$ws.Range("M134").Value = -10532.4
$ws.Range("N134").Value = -22426.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H42").Value = 30000
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 30000
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("N42").Value = -31186
$ws.Range("L42").Value = 30000
$ws.Range("H132").Value = 5684.467
$ws.Range("I132").Value = 4591.609
$ws.Range("K132").Value = 13774.827
$ws.Range("M132").Value = -11244.827
$ws.Range("H133").Value = 58680.727
$ws.Range("J133").Value = 60061
$ws.Range("L133").Value = 60061
$ws.Range("N133").Value = -65121
$ws.Range("H134").Value = 2947.9167
$ws.Range("I134").Value = 2282.6667
$ws.Range("J134").Value = 7604.6665
$ws.Range("K134").Value = 6848.000100000001
$ws.Range("L134").Value = 22813.9995
$ws.Range("M134").Value = -4313.000100000001
$ws.Range("N134").Value = -27883.9995
$ws.Range("H141").Value = 277166.38
$ws.Range("J141").Value = 299883
$ws.Range("L141").Value = 299883
$ws.Range("N141").Value = -310243

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 7217.364
$ws.Range("I56").Value = 7217.364
$ws.Range("K56").Value = 7217.364
$ws.Range("M56").Value = -6687.364
$ws.Range("H113").Value = 3220.1738
$ws.Range("J113").Value = 3252.9092
$ws.Range("L113").Value = 9758.7276
$ws.Range("N113").Value = -14098.7276
$ws.Range("H117").Value = 2201.2
$ws.Range("I117").Value = 243.5
$ws.Range("K117").Value = 730.5
$ws.Range("M117").Value = 2711.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5177.9165
$ws.Range("I122").Value = 4680.381
$ws.Range("J122").Value = 8660.666999999999
$ws.Range("K122").Value = 14041.143
$ws.Range("L122").Value = 25982.001
$ws.Range("M122").Value = -11591.143
$ws.Range("N122").Value = -30882.001
$ws.Range("H126").Value = 4441.7915
$ws.Range("I126").Value = 3059.7
$ws.Range("K126").Value = 9179.099999999999
$ws.Range("M126").Value = -6709.099999999999
$ws.Range("H132").Value = 5582.6
$ws.Range("I132").Value = 2224.75
$ws.Range("J132").Value = 19014
$ws.Range("K132").Value = 6674.25
$ws.Range("L132").Value = 57042
$ws.Range("M132").Value = -4144.25
$ws.Range("N132").Value = -62102
$ws.Range("H138").Value = 68656.664
$ws.Range("J138").Value = 68656.664
$ws.Range("L138").Value = 68656.664
$ws.Range("N138").Value = -78936.664
$ws.Range("H141").Value = 59960.6
$ws.Range("J141").Value = 59960.6
$ws.Range("L141").Value = 59960.6
$ws.Range("N141").Value = -70320.60000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 16972
$ws.Range("I22").Value = 1416.1666
$ws.Range("K22").Value = 1416.1666
$ws.Range("M22").Value = -1121.1666
$ws.Range("H27").Value = 16972
$ws.Range("I27").Value = 1416.1666
$ws.Range("K27").Value = 1416.1666
$ws.Range("M27").Value = -1309.1666
$ws.Range("H46").Value = 4452.8184
$ws.Range("I46").Value = 2330
$ws.Range("J46").Value = 5248.875
$ws.Range("K46").Value = 2330
$ws.Range("L46").Value = 5248.875
$ws.Range("M46").Value = -2142
$ws.Range("N46").Value = -5624.875
$ws.Range("H55").Value = 2466.8
$ws.Range("I55").Value = 1465
$ws.Range("K55").Value = 1465
$ws.Range("M55").Value = -1292
$ws.Range("H69").Value = 30999.5
$ws.Range("J69").Value = 25000
$ws.Range("L69").Value = 25000
$ws.Range("N69").Value = -26622
$ws.Range("H72").Value = 30999.5
$ws.Range("J72").Value = 25000
$ws.Range("L72").Value = 75000
$ws.Range("N72").Value = -83112
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("L118").Value = 0
$ws.Range("H122").Value = 166296.84
$ws.Range("I122").Value = 228008.83
$ws.Range("K122").Value = 684026.49
$ws.Range("M122").Value = -681576.49
$ws.Range("H132").Value = 8358.666999999999
$ws.Range("I132").Value = 5787.5
$ws.Range("K132").Value = 17362.5
$ws.Range("M132").Value = -14832.5
$ws.Range("H136").Value = 5231.885
$ws.Range("I136").Value = 3819.7646
$ws.Range("J136").Value = 7899.222
$ws.Range("K136").Value = 11459.2938
$ws.Range("L136").Value = 23697.666
$ws.Range("M136").Value = -8909.293799999999
$ws.Range("N136").Value = -28797.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 14002.5
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("H141").Value = 99685.8
$ws.Range("J141").Value = 111357.25
$ws.Range("L141").Value = 111357.25
$ws.Range("N141").Value = -121717.25
